$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("APIData")

# Clear the existing used range (old data had columns A:D, rows 1:2)
$ws.Cells.Clear()

# Write new header row and data row, ensuring shared-string insertion
# order matches TestCase, TempVariance, Temperature(in Celsius)
$ws.Range("A1").Value = "TestCase"
$ws.Range("A2").Value = "TempVariance"
$ws.Range("B1").Value = "Temperature(in Celsius)"
$ws.Range("B2").Value = 29.26

# Set column widths to match target layout (closest achievable values,
# since column width is stored/rounded to the nearest pixel internally)
$ws.Columns.Item(1).ColumnWidth = 10.666666666666666
$ws.Columns.Item(2).ColumnWidth = 18.833333333333336
$ws.Columns.Item(3).ColumnWidth = 16.833333333333336

# Update selection to B2 as in the target worksheet view
$ws.Activate()
$ws.Range("B2").Select()

$wb.Save()
